# Apply the "Shipping Logic still in progress - Slow Grind - 24 Dec 2025" update:
#  1) Rewrite the Settings sheet (A2:B17) with the new parameter list/values.
#  2) Update the Log sheet timestamps (A2:A33) and a few Details cells
#     (C4, C16, C25) to reflect the new run.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Settings sheet
# ---------------------------------------------------------------------------
$settings = $wb.Worksheets.Item("Settings")

$settings.Range("A2").Value = "horizon_days"
$settings.Range("B2").Value = 365

$settings.Range("A3").Value = "random_opening"
$settings.Range("B3").Value = $true

$settings.Range("A4").Value = "random_seed"
$settings.Range("B4").Value = ""

$settings.Range("A5").Value = "progress_step_pct"
$settings.Range("B5").Value = 5

$settings.Range("A6").Value = "demand_truck_load_tons"
$settings.Range("B6").Value = 25

$settings.Range("A7").Value = "demand_step_hours"
$settings.Range("B7").Value = 1

$settings.Range("A8").Value = "require_full_payload"
$settings.Range("B8").Value = $true

$settings.Range("A9").Value = "ship_idle_wait_h"
$settings.Range("B9").Value = 1

$settings.Range("A10").Value = "ship_max_wait_product_h"
$settings.Range("B10").Value = 24

$settings.Range("A11").Value = "transporter_wait_h"
$settings.Range("B11").Value = 1

$settings.Range("A12").Value = "mean_breakdown_duration"
$settings.Range("B12").Value = 3

$settings.Range("A13").Value = "make_output_choice"
$settings.Range("B13").Value = "min_fill_pct"

$settings.Range("A14").Value = "step_hours"
$settings.Range("B14").Value = 1

$settings.Range("A15").Value = "write_plots"
$settings.Range("B15").Value = $true

$settings.Range("A16").Value = "write_csvs"
$settings.Range("B16").Value = $true

$settings.Range("A17").Value = "out_dir"
$settings.Range("B17").Value = "sim_outputs"

# ---------------------------------------------------------------------------
# 2) Log sheet
# ---------------------------------------------------------------------------
$log = $wb.Worksheets.Item("Log")

$log.Range("A2").Value = "2025-12-23 22:54:58"
$log.Range("A3").Value = "2025-12-23 22:54:58"
$log.Range("A4").Value = "2025-12-23 22:54:58"
$log.Range("A5").Value = "2025-12-23 22:54:58"
$log.Range("A6").Value = "2025-12-23 22:54:58"
$log.Range("A7").Value = "2025-12-23 22:54:58"
$log.Range("A8").Value = "2025-12-23 22:54:58"
$log.Range("A9").Value = "2025-12-23 22:54:58"
$log.Range("A10").Value = "2025-12-23 22:54:58"
$log.Range("A11").Value = "2025-12-23 22:54:59"
$log.Range("A12").Value = "2025-12-23 22:54:59"
$log.Range("A13").Value = "2025-12-23 22:54:59"
$log.Range("A14").Value = "2025-12-23 22:54:59"
$log.Range("A15").Value = "2025-12-23 22:54:59"
$log.Range("A16").Value = "2025-12-23 22:54:59"
$log.Range("A17").Value = "2025-12-23 22:54:59"
$log.Range("A18").Value = "2025-12-23 22:54:59"
$log.Range("A19").Value = "2025-12-23 22:54:59"
$log.Range("A20").Value = "2025-12-23 22:54:59"
$log.Range("A21").Value = "2025-12-23 22:54:59"
$log.Range("A22").Value = "2025-12-23 22:54:59"
$log.Range("A23").Value = "2025-12-23 22:54:59"
$log.Range("A24").Value = "2025-12-23 22:54:59"
$log.Range("A25").Value = "2025-12-23 22:54:59"
$log.Range("A26").Value = "2025-12-23 22:55:00"
$log.Range("A27").Value = "2025-12-23 22:55:00"
$log.Range("A28").Value = "2025-12-23 22:55:00"
$log.Range("A29").Value = "2025-12-23 22:55:01"
$log.Range("A30").Value = "2025-12-23 22:55:01"
$log.Range("A31").Value = "2025-12-23 22:55:01"
$log.Range("A32").Value = "2025-12-23 22:55:02"
$log.Range("A33").Value = "2025-12-23 22:55:02"

$log.Range("C4").Value = "Settings -> Settings, rows=25"
$log.Range("C16").Value = "Settings rows=25 -> generated_model_inputs.xlsx"
$log.Range("C25").Value = "added=16, updated=0"
